# Generate Report for handoff
# The "eec112a9-8a1c-41d7-aad5-738d1b829fd7.md" source file has been handed
# off again: its status moves from "Handed back: in sync with en-US" to
# "Ready for handoff" on every language sheet, and the zh-cn / de-de sheets
# record the new handoff file + handoff datetime for that row.

$wb = $excel.ActiveWorkbook

$statusReadyForHandoff = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: row for eec112a9-8a1c-41d7-aad5-738d1b829fd7.md (row 3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusReadyForHandoff
$wsOverview.Range("C3").Value = $statusReadyForHandoff

# ---------------------------------------------------------------------
# zh-cn sheet: row for eec112a9-8a1c-41d7-aad5-738d1b829fd7.md (row 3)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $statusReadyForHandoff
$wsZhCn.Range("C3").Value = "eec112a9-8a1c-41d7-aad5-738d1b829fd7.4610716fb7e1428f87d9ecfc5d57e86846380d74.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-01-14 02:37:22"

# Rebuild the hyperlinks for this sheet so the C3 hyperlink's display text
# reflects the new handoff file while every other hyperlink (address + rId
# order) stays exactly as it was.
$wsZhCn.Range("A1:I4").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3141594a8ee7bdcf32b7daccf45d9fcfdcfd497d/e2e/42db75d7-9895-4a12-b95c-d3a35481d143.md", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfc07c4fd918b0b1cdb4a4999170fa034958bcf6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7c3e974cafe9b7e0440941c6272825d2a7f70b72/e2e/42db75d7-9895-4a12-b95c-d3a35481d143.md", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0331053de06f1e5d496b2d4cb42b0e5c514d00e4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3141594a8ee7bdcf32b7daccf45d9fcfdcfd497d/e2e/eec112a9-8a1c-41d7-aad5-738d1b829fd7.md", "", "", "eec112a9-8a1c-41d7-aad5-738d1b829fd7.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dfc07c4fd918b0b1cdb4a4999170fa034958bcf6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf", "", "", "eec112a9-8a1c-41d7-aad5-738d1b829fd7.4610716fb7e1428f87d9ecfc5d57e86846380d74.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/7c3e974cafe9b7e0440941c6272825d2a7f70b72/e2e/42db75d7-9895-4a12-b95c-d3a35481d143.md", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0331053de06f1e5d496b2d4cb42b0e5c514d00e4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.zh-cn.xlf")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3141594a8ee7bdcf32b7daccf45d9fcfdcfd497d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet: row for eec112a9-8a1c-41d7-aad5-738d1b829fd7.md (row 3)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $statusReadyForHandoff
$wsDeDe.Range("C3").Value = "eec112a9-8a1c-41d7-aad5-738d1b829fd7.4610716fb7e1428f87d9ecfc5d57e86846380d74.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-01-14 02:37:35"

$wsDeDe.Range("A1:I4").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/3141594a8ee7bdcf32b7daccf45d9fcfdcfd497d/e2e/42db75d7-9895-4a12-b95c-d3a35481d143.md", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fadb95c3876c02cbbfb88eca17a0e00fe528b79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fdc5e40528844d381ca5693e874cc81045ad43ec/e2e/42db75d7-9895-4a12-b95c-d3a35481d143.md", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/233d27a7013d98fd83ccd856466183cec588b858/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/3141594a8ee7bdcf32b7daccf45d9fcfdcfd497d/e2e/eec112a9-8a1c-41d7-aad5-738d1b829fd7.md", "", "", "eec112a9-8a1c-41d7-aad5-738d1b829fd7.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5fadb95c3876c02cbbfb88eca17a0e00fe528b79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf", "", "", "eec112a9-8a1c-41d7-aad5-738d1b829fd7.4610716fb7e1428f87d9ecfc5d57e86846380d74.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fdc5e40528844d381ca5693e874cc81045ad43ec/e2e/42db75d7-9895-4a12-b95c-d3a35481d143.md", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/233d27a7013d98fd83ccd856466183cec588b858/ol-handback/OpenLocalizationTestOrg/oltest.de-de/yuwzho/42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf", "", "", "42db75d7-9895-4a12-b95c-d3a35481d143.610c8d1acb69e0aab4c42f9bd02b9cb2816d5378.de-de.xlf")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/3141594a8ee7bdcf32b7daccf45d9fcfdcfd497d/.localization-config", "", "", ".localization-config")
